# Insert a new daily-price row for "Plátano" (Vega Monumental Concepción)
# right above the current row 288. This shifts the existing rows 288-402
# down to 289-403 (carrying all of their original values with them) and
# the sheet's used range grows from A1:T402 to A1:T403.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(288).Insert()

# Populate the newly inserted row 288 with the new observation.
$ws.Cells.Item(288, 1).Value  = 11
$ws.Cells.Item(288, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(288, 3).Value  = "Bíobío"
$ws.Cells.Item(288, 4).Value  = 44553
$ws.Cells.Item(288, 5).Value  = 8
$ws.Cells.Item(288, 6).Value  = "Fruta"
$ws.Cells.Item(288, 7).Value  = 100108
$ws.Cells.Item(288, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(288, 9).Value  = 100108006
$ws.Cells.Item(288, 10).Value = "Plátano"
$ws.Cells.Item(288, 11).Value = "Sin especificar"
$ws.Cells.Item(288, 12).Value = "Pintón"
$ws.Cells.Item(288, 13).Value = 750
$ws.Cells.Item(288, 14).Value = 10000
$ws.Cells.Item(288, 15).Value = 11000
$ws.Cells.Item(288, 16).Value = 10467
$ws.Cells.Item(288, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(288, 18).Value = "Ecuador"
$ws.Cells.Item(288, 19).Value = 523
$ws.Cells.Item(288, 20).Value = 20
